$wb = $excel.ActiveWorkbook

# --- Update the Metrics sheet values ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 320328.74
$metrics.Range("B3").Value = 262495.15000000002
$metrics.Range("B4").Value = 102106.30999999998
$metrics.Range("B5").Value = 12788
$metrics.Range("B6").Value = 4687460.21
$metrics.Range("B7").Value = 3952313.8199999994
$metrics.Range("B8").Value = 1372708.45
$metrics.Range("B9").Value = 181789
$metrics.Range("B10").Value = 33152784.010999821
$metrics.Range("B11").Value = 31227535.34
$metrics.Range("B12").Value = 11654417.340000002
$metrics.Range("B13").Value = 1279416

# Force recalculation so dependent formulas on other sheets (e.g. "today")
# pick up the new Metrics values.
$excel.CalculateFullRebuild()

# --- Update the selection (active cell) on each sheet, matching the diff ---
[void]$metrics.Range("D39").Select()

$today = $wb.Worksheets.Item("today")
[void]$today.Activate()
[void]$today.Range("E5").Select()
